# Home Screen test cases are added
#
# This script mutates the "Profile.xlsx" test-case workbook:
#  - Sheet "MA_AccountEdit1": the login credentials used by the test are
#    swapped out (old webapps.com test account -> nfhslearn.com admin
#    account) and a "Results" value is recorded.
#  - Sheet "Test Cases": a "Pass" result is recorded for the test row.
#  - The active sheet/selection bookkeeping is updated so that "Test Cases"
#    (not "MA_AccountEdit1") is the tab shown when the workbook re-opens.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Test Cases"
$ws2 = $wb.Worksheets.Item(2)   # "MA_AccountEdit1"

# ---------------------------------------------------------------------
# Sheet "MA_AccountEdit1" (2nd sheet)
# ---------------------------------------------------------------------

# Username (A2): swap the old webapps.com address for the nfhslearn.com
# admin account, and repoint its mailto hyperlink to match.
$ws2.Range("A2").Value = "admin@nfhslearn.com"
$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:admin@nfhslearn.com")
$ws2.Range("A2").Style = "Hyperlink"

# Password (B2): swap in the new password, and give it a mailto
# hyperlink (matching the look of the username cell next to it).
$ws2.Range("B2").Value = "nfhslearn@6186"
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:nfhslearn@6186")
$ws2.Range("B2").Style = "Hyperlink"

# Results (D2): record the outcome of the run.
$ws2.Range("D2").Value = "pass"
$ws2.Range("D2").Style = "Normal"

# Widen the Password column a bit to fit the new value.
$ws2.Columns.Item(2).ColumnWidth = 30.02

# ---------------------------------------------------------------------
# Sheet "Test Cases" (1st sheet)
# ---------------------------------------------------------------------

# Results (E2): record the outcome of the run.
$ws1.Range("E2").Value = "Pass"

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping
#   Select "MA_AccountEdit1" first (remembering its own cell selection),
#   then finish on "Test Cases" so that it ends up the active/visible tab.
# ---------------------------------------------------------------------
$ws2.Select()
$ws2.Range("C5").Select()

$ws1.Select()
$ws1.Range("B4").Select()
